# Applies the "added webtable handle concept and javascript executor concept" edit:
#  - Removes the now-unused "Status" column (K:N) from RegTestData.
#  - Adds a second worksheet "TableData" holding a small web-table style
#    dataset (Company Name / Contact Name), mirroring a typical Selenium
#    WebTable-handling example fixture.
#  - Leaves the active selection on RegTestData at I11 and activates the
#    new TableData sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- RegTestData: drop the trailing Status column (K) and the three blank
#     formatted-but-empty columns after it (L, M, N) ---
$ws1.Range("K1:N1").EntireColumn.Delete()

# Update the remembered selection on RegTestData to I11 (matches target).
$ws1.Range("I11").Select() | Out-Null

# --- Add the new TableData worksheet right after RegTestData ---
$active = $wb.ActiveSheet
$ws2 = $wb.Worksheets.Add($null, $active)
$ws2.Name = "TableData"

# Header row
$ws2.Range("A1").Value = "Company Name"
$ws2.Range("B1").Value = "Contact Name"

# Data rows (classic "Customers" sample set used for WebTable examples)
$data = @(
    @("Alfreds Futterkiste", "Maria Anders"),
    @("Centro comercial Moctezuma", "Francisco Chang"),
    @("Ernst Handel", "Roland Mendel"),
    @("Island Trading", "Helen Bennett"),
    @("Laughing Bacchus Winecellars", "Yoshi Tannamuri"),
    @("Magazzini Alimentari Riuniti", "Giovanni Rovelli")
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# TableData becomes the active/visible sheet.
$ws2.Activate() | Out-Null
$ws2.Range("A1").Select() | Out-Null
